$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values can look like plain numbers (e.g. "336.26"); force the
# Price column to a text format before writing so Excel does not
# silently convert them to numeric cells, then restore the default
# "Normal" style so no residual number-format style is left behind.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.251.20"
$ws.Range("E2").Value = "  +1.33%  "

$ws.Range("D3").Value = "1.805.84"
$ws.Range("E3").Value = "  +3.02%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").Value = "336.26"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("D6").Value = "0.9989"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").Value = "0.4618"
$ws.Range("E7").Value = "  +20.96%  "

$ws.Range("D8").Value = "0.3707"
$ws.Range("E8").Value = "  +9.38%  "

$ws.Range("D9").Value = "45.15"
$ws.Range("E9").Value = "  +1.14%  "

$ws.Range("D10").Value = "1.153"
$ws.Range("E10").Value = "  +3.70%  "

$ws.Range("D11").Value = "0.07643"
$ws.Range("E11").Value = "  +5.97%  "

$ws.Range("D12").Value = "22.43"
$ws.Range("E12").Value = "  -0.21%  "

$ws.Range("D13").Value = "1.001"
$ws.Range("E13").Value = "  -0.24%  "

$ws.Range("D14").Value = "6.344"
$ws.Range("E14").Value = "  +3.08%  "

$ws.Range("D15").Value = "7.482"
$ws.Range("E15").Value = "  +4.86%  "

$ws.Range("D16").Value = "1.805.76"
$ws.Range("E16").Value = "  +2.81%  "

$ws.Range("E17").Value = "  +3.98%  "

$ws.Range("D18").Value = "0.06722"
$ws.Range("E18").Value = "  +1.79%  "

$ws.Range("D19").Value = "82.00"
$ws.Range("E19").Value = "  +3.88%  "

$ws.Range("D20").Value = "0.9989"
$ws.Range("E20").Value = "  -0.20%  "

$ws.Range("D21").Value = "17.50"
$ws.Range("E21").Value = "  +4.96%  "

$ws.Range("D22").Value = "6.426"
$ws.Range("E22").Value = "  +3.26%  "

$ws.Range("D23").Value = "28.244.17"
$ws.Range("E23").Value = "  +1.24%  "

$ws.Range("D24").Value = "11.89"
$ws.Range("E24").Value = "  +2.31%  "

$ws.Range("D25").Value = "2.409"
$ws.Range("E25").Value = "  +0.98%  "

$ws.Range("D26").Value = "20.81"
$ws.Range("E26").Value = "  +5.04%  "

$ws.Range("D27").Value = "153.67"
$ws.Range("E27").Value = "  +0.97%  "

$ws.Range("D28").Value = "2.383"
$ws.Range("E28").Value = "  +2.95%  "

$ws.Range("D29").Value = "2.011.35"
$ws.Range("E29").Value = "  +2.72%  "

$ws.Range("D30").Value = "133.49"
$ws.Range("E30").Value = "  +1.22%  "

$ws.Range("D31").Value = "1.259"
$ws.Range("E31").Value = "  -0.84%  "

$ws.Range("D32").Value = "4.034"
$ws.Range("E32").Value = "  +0.31%  "

$ws.Range("D33").Value = "0.09574"
$ws.Range("E33").Value = "  +8.82%  "

$ws.Range("D34").Value = "5.871"
$ws.Range("E34").Value = "  +0.78%  "

$ws.Range("E35").Value = "  +5.83%  "

$ws.Range("D36").Value = "12.14"
$ws.Range("E36").Value = "  -0.47%  "

$ws.Range("D37").Value = "0.06376"
$ws.Range("E37").Value = "  +3.63%  "

$ws.Range("D38").Value = "0.02358"
$ws.Range("E38").Value = "  +3.33%  "

$ws.Range("D39").Value = "5.270"
$ws.Range("E39").Value = "  +2.86%  "

$ws.Range("D40").Value = "0.6657"
$ws.Range("E40").Value = "  +1.29%  "

$ws.Range("D41").Value = "1.239"
$ws.Range("E41").Value = "  +2.53%  "

$ws.Range("D42").Value = "1.497"
$ws.Range("E42").Value = "  -0.34%  "

$ws.Range("D43").Value = "8.181"
$ws.Range("E43").Value = "  +2.44%  "

$ws.Range("D44").Value = "14.27"
$ws.Range("E44").Value = "  +4.72%  "

$ws.Range("E45").Value = "  -0.28%  "

$ws.Range("D46").Value = "0.6143"
$ws.Range("E46").Value = "  +1.44%  "

$ws.Range("D47").Value = "3.827"
$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("D48").Value = "130.00"
$ws.Range("E48").Value = "  +3.07%  "

$ws.Range("D49").Value = "2.051"
$ws.Range("E49").Value = "  +2.60%  "

$ws.Range("D50").Value = "0.07161"
$ws.Range("E50").Value = "  +2.59%  "

$ws.Range("E51").Value = "  +0.57%  "

$priceRange.Style = "Normal"
